# edit.ps1
#
# "Update countries & provincias Spain" - refresh of the COVID-19 country table.
#
# The source feed was re-pulled later in the day, which:
#   1. Bumped the "Datos actualizados..." timestamp (cell A1).
#   2. Re-ranked several countries by total cases, so some rows now show a
#      different country name (column A) than before - "Bolivia"/"Ecuador",
#      "Bahamas"/"Nueva Zelanda" and "Birmania"/"Burundi" each swap with their
#      neighbouring row, while "Islas Virgenes Britanicas" jumps up several
#      places, pushing the rows below it down by one.
#   3. Refreshed the statistics columns (B:H = Casos totales, Nuevos casos,
#      Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) for
#      every row whose figures changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Timestamp -------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 22 de Agosto de 2020 a las 05:00"

# --- 2. Country label swaps caused by the re-ranking --------------------
$countrySwaps = [ordered]@{
    "A29" = "Bolivia"
    "A30" = "Ecuador"
    "A142" = "Bahamas"
    "A143" = "Nueva Zelanda"
    "A171" = "Birmania"
    "A172" = "Burundi"
    "A207" = "Islas Virgenes Britanicas"
    "A208" = "Dominica"
    "A209" = "Islas Virgenes de los Estados Unidos"
    "A210" = "San Cristobal y Nieves"
    "A211" = "San Bartolome"
    "A212" = "Groenlandia"
    "A213" = "Bonaire, San Eustaquio y Saba"
    "A214" = "Islas Malvinas"
    "A215" = "Montserrat"
}
foreach ($cellRef in $countrySwaps.Keys) {
    $ws.Range($cellRef).Value = $countrySwaps[$cellRef]
}

# --- 3. Refreshed statistics ---------------------------------------------
$statUpdates = [ordered]@{
    23 = @{ D=208950; E=14743 }
    29 = @{ B=107435; C=1370; D=42141; E=60928; G=61; H=4366 }
    30 = @{ B=106481; D=87730; E=12503; H=6248 }
    31 = @{ B=104313; C=242; E=14978 }
    40 = @{ B=80894; C=716; D=18165; E=52744; G=9; H=9985 }
    51 = @{ B=53381; C=562; D=8271; E=43478; G=13; H=1632 }
    86 = @{ B=12536; D=7007; E=5347; H=182 }
    142 = @{ B=1703; C=0; D=221; E=1455; H=27 }
    143 = @{ B=1671; C=6; D=1538; E=111; H=22 }
    171 = @{ B=435; C=16; D=337; E=92; H=6 }
    172 = @{ B=426; D=336; E=89; H=1 }
    184 = @{ D=263; E=10 }
    207 = @{ B=21; C=9; D=8; E=12; H=1 }
    208 = @{ B=18; D=18; E=0 }
    209 = @{ D=0; E=17 }
    210 = @{ B=17; D=17; E=0 }
    211 = @{ B=16; D=9; E=7 }
    212 = @{ B=14; D=14; E=0 }
    213 = @{ D=7; E=6 }
    214 = @{ D=13; H=0 }
    215 = @{ B=13; D=12; E=0 }
}
foreach ($row in $statUpdates.Keys) {
    $cols = $statUpdates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
